$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Kommunal)
$ws.Range("B2").Value = 0.0331
$ws.Range("C2").Value = 2.1583
$ws.Range("D2").Value = 152.2319
$ws.Range("E2").Value = 0.9397
$ws.Range("F2").Value = 155.363

# Row 3 (Privat)
$ws.Range("B3").Value = 0.4055
$ws.Range("C3").Value = 26.4743
$ws.Range("D3").Value = 413.7398
$ws.Range("E3").Value = 3.7279
$ws.Range("F3").Value = 444.3475
